# Apply the data updates from the "Sales Jan-June" / "Sales July-Dec" table
# on Sheet1 of the CFTemplate workbook. Only the actual numeric values
# changed between the two revisions (everything else in the recorded diff
# is boilerplate that Excel itself rewrites on save, e.g. namespace/version
# bumps, style-serialisation order, calcPr fullPrecision flag, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value  = 60000
$ws.Range("D10").Value = 32000
$ws.Range("D17").Value = 58500
$ws.Range("D23").Value = 65700
$ws.Range("D36").Value = 51500
$ws.Range("D42").Value = 58300
$ws.Range("D43").Value = 62400
